$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-08-14 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-15 Thursday", 2) | Out-Null

# Update each table cell precisely by row/column to avoid ambiguity
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text = "66-25=41"
$table.Cell(1, 2).Range.Text = "11+39=50"
$table.Cell(1, 3).Range.Text = "46+21=67"
$table.Cell(1, 4).Range.Text = "58+12=70"
$table.Cell(1, 5).Range.Text = "29+10=39"
$table.Cell(2, 1).Range.Text = "62-2=60"
$table.Cell(2, 2).Range.Text = "41+38=79"
$table.Cell(2, 3).Range.Text = "94-44=50"
$table.Cell(2, 4).Range.Text = "60-24=36"
$table.Cell(2, 5).Range.Text = "12-7=5"
$table.Cell(3, 1).Range.Text = "99-19=80"
$table.Cell(3, 2).Range.Text = "32+65=97"
$table.Cell(3, 3).Range.Text = "75-1=74"
$table.Cell(3, 4).Range.Text = "26+71=97"
$table.Cell(3, 5).Range.Text = "12+41=53"
$table.Cell(4, 1).Range.Text = "62-61=1"
$table.Cell(4, 2).Range.Text = "38-8=30"
$table.Cell(4, 3).Range.Text = "42-6=36"
$table.Cell(4, 4).Range.Text = "86-1=85"
$table.Cell(4, 5).Range.Text = "59+25=84"
$table.Cell(5, 1).Range.Text = "53-30=23"
$table.Cell(5, 2).Range.Text = "48-10=38"
$table.Cell(5, 3).Range.Text = "35-30=5"
$table.Cell(5, 4).Range.Text = "60-18=42"
$table.Cell(5, 5).Range.Text = "44-18=26"
$table.Cell(6, 1).Range.Text = "17+30=47"
$table.Cell(6, 2).Range.Text = "1+40=41"
$table.Cell(6, 3).Range.Text = "83-1=82"
$table.Cell(6, 4).Range.Text = "67-18=49"
$table.Cell(6, 5).Range.Text = "22+21=43"
$table.Cell(7, 1).Range.Text = "9+27=36"
$table.Cell(7, 2).Range.Text = "71-65=6"
$table.Cell(7, 3).Range.Text = "92-82=10"
$table.Cell(7, 4).Range.Text = "29+36=65"
$table.Cell(7, 5).Range.Text = "94-13=81"
$table.Cell(8, 1).Range.Text = "88-65=23"
$table.Cell(8, 2).Range.Text = "25+32=57"
$table.Cell(8, 3).Range.Text = "60+31=91"
$table.Cell(8, 4).Range.Text = "2+71=73"
$table.Cell(8, 5).Range.Text = "61-60=1"
$table.Cell(9, 1).Range.Text = "24+42=66"
$table.Cell(9, 2).Range.Text = "16+69=85"
$table.Cell(9, 3).Range.Text = "75-51=24"
$table.Cell(9, 4).Range.Text = "83-56=27"
$table.Cell(9, 5).Range.Text = "21+48=69"
$table.Cell(10, 1).Range.Text = "9+24=33"
$table.Cell(10, 2).Range.Text = "71-57=14"
$table.Cell(10, 3).Range.Text = "42+37=79"
$table.Cell(10, 4).Range.Text = "62-29=33"
$table.Cell(10, 5).Range.Text = "46-46=0"
$table.Cell(11, 1).Range.Text = "43-3=40"
$table.Cell(11, 2).Range.Text = "46+50=96"
$table.Cell(11, 3).Range.Text = "71+4=75"
$table.Cell(11, 4).Range.Text = "53-36=17"
$table.Cell(11, 5).Range.Text = "50-40=10"
$table.Cell(12, 1).Range.Text = "12+50=62"
$table.Cell(12, 2).Range.Text = "84+7=91"
$table.Cell(12, 3).Range.Text = "40-23=17"
$table.Cell(12, 4).Range.Text = "46-3=43"
$table.Cell(12, 5).Range.Text = "89-63=26"
$table.Cell(13, 1).Range.Text = "21+50=71"
$table.Cell(13, 2).Range.Text = "81-14=67"
$table.Cell(13, 3).Range.Text = "47-40=7"
$table.Cell(13, 4).Range.Text = "75-11=64"
$table.Cell(13, 5).Range.Text = "56+28=84"
$table.Cell(14, 1).Range.Text = "46-18=28"
$table.Cell(14, 2).Range.Text = "49+50=99"
$table.Cell(14, 3).Range.Text = "35+17=52"
$table.Cell(14, 4).Range.Text = "93-67=26"
$table.Cell(14, 5).Range.Text = "27+6=33"
$table.Cell(15, 1).Range.Text = "98-30=68"
$table.Cell(15, 2).Range.Text = "61-46=15"
$table.Cell(15, 3).Range.Text = "26+40=66"
$table.Cell(15, 4).Range.Text = "40+17=57"
$table.Cell(15, 5).Range.Text = "29-24=5"
$table.Cell(16, 1).Range.Text = "47+39=86"
$table.Cell(16, 2).Range.Text = "90-63=27"
$table.Cell(16, 3).Range.Text = "11+54=65"
$table.Cell(16, 4).Range.Text = "37+60=97"
$table.Cell(16, 5).Range.Text = "8+79=87"
$table.Cell(17, 1).Range.Text = "81-68=13"
$table.Cell(17, 2).Range.Text = "13+49=62"
$table.Cell(17, 3).Range.Text = "8+20=28"
$table.Cell(17, 4).Range.Text = "2+38=40"
$table.Cell(17, 5).Range.Text = "14-13=1"
$table.Cell(18, 1).Range.Text = "33-15=18"
$table.Cell(18, 2).Range.Text = "79+7=86"
$table.Cell(18, 3).Range.Text = "82-82=0"
$table.Cell(18, 4).Range.Text = "52+26=78"
$table.Cell(18, 5).Range.Text = "70+22=92"
$table.Cell(19, 1).Range.Text = "69-3=66"
$table.Cell(19, 2).Range.Text = "60-6=54"
$table.Cell(19, 3).Range.Text = "69-65=4"
$table.Cell(19, 4).Range.Text = "30+22=52"
$table.Cell(19, 5).Range.Text = "55-16=39"
$table.Cell(20, 1).Range.Text = "27+44=71"
$table.Cell(20, 2).Range.Text = "59-54=5"
$table.Cell(20, 3).Range.Text = "40+41=81"
$table.Cell(20, 4).Range.Text = "57-26=31"
$table.Cell(20, 5).Range.Text = "3+42=45"
